$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the revised GDP figures for existing years (rows 12-30)
$ws.Range("B12").Value = 401407.69300000003
$ws.Range("B13").Value = 411155.21100000001
$ws.Range("B14").Value = 425385.58600000001
$ws.Range("B15").Value = 449004.15600000002
$ws.Range("B16").Value = 472993.55200000003
$ws.Range("B17").Value = 500969.50199999998
$ws.Range("B18").Value = 523165.163
$ws.Range("B19").Value = 520597.01500000001
$ws.Range("B20").Value = 510758.95799999998
$ws.Range("B21").Value = 523098.04100000003
$ws.Range("B22").Value = 542384.78
$ws.Range("B23").Value = 570954.87399999995
$ws.Range("B24").Value = 583070.71999999997
$ws.Range("B25").Value = 609478.18000000005
$ws.Range("B26").Value = 636938.86499999999
$ws.Range("B27").Value = 647455.652
$ws.Range("B28").Value = 667153.495
$ws.Range("B29").Value = 701954.96799999999
$ws.Range("B30").Value = 721906.95

# Add the new observation for 2020-01-01
$ws.Range("A31").Value = 43831
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B31").Value = 692988.23699999996
$ws.Range("B31").NumberFormat = "0.000"

# Update the sheet selection state
$ws.Range("A1:B1048576").Select()
